# Auto-generated: applies the Diabolos_Profits data refresh described in the commit diff.
# Each leve's currentAveragePrice / NQ / HQ price+profit columns are refreshed with
# newer market-board snapshot values (static data, no formulas in this workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 64206.5
$ws.Range("J112").Value = 1820.2667
$ws.Range("L112").Value = 5460.800099999999
$ws.Range("N112").Value = -7676.800099999999
# Row 129
$ws.Range("H129").Value = 2013.5555
$ws.Range("I129").Value = 1542.3334
$ws.Range("K129").Value = 4627.0002
$ws.Range("M129").Value = 372.9997999999996
# Row 137
$ws.Range("H137").Value = 5361.8184
$ws.Range("I137").Value = 4198.3076
$ws.Range("J137").Value = 7042.4443
$ws.Range("K137").Value = 12594.9228
$ws.Range("L137").Value = 21127.3329
$ws.Range("M137").Value = -10044.9228
$ws.Range("N137").Value = -26227.3329

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1103.5385
$ws.Range("I2").Value = 445.5
$ws.Range("J2").Value = 9000
$ws.Range("K2").Value = 445.5
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = -332.5
$ws.Range("N2").Value = -9226
# Row 32
$ws.Range("H32").Value = 17274.793
$ws.Range("I32").Value = 12628.556
$ws.Range("K32").Value = 12628.556
$ws.Range("M32").Value = -12341.556
# Row 45
$ws.Range("H45").Value = 266226.72
$ws.Range("I45").Value = 310331.22
$ws.Range("K45").Value = 310331.22
$ws.Range("M45").Value = -309954.22
# Row 74
$ws.Range("H74").Value = 3070.9355
$ws.Range("I74").Value = 2346.8823
$ws.Range("J74").Value = 3950.1428
$ws.Range("K74").Value = 2346.8823
$ws.Range("L74").Value = 3950.1428
$ws.Range("M74").Value = -1472.8823
$ws.Range("N74").Value = -5698.1428
# Row 77
$ws.Range("H77").Value = 3070.9355
$ws.Range("I77").Value = 2346.8823
$ws.Range("J77").Value = 3950.1428
$ws.Range("K77").Value = 11734.4115
$ws.Range("L77").Value = 19750.714
$ws.Range("M77").Value = -7366.411500000002
$ws.Range("N77").Value = -28486.714
# Row 110
$ws.Range("H110").Value = 34484750
$ws.Range("I110").Value = 47620412
$ws.Range("K110").Value = 47620412
$ws.Range("M110").Value = -47618367
# Row 116
$ws.Range("H116").Value = 1103.5385
$ws.Range("I116").Value = 445.5
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 445.5
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = 1848.5
$ws.Range("N116").Value = -13588
# Row 139
$ws.Range("H139").Value = 66238.664
$ws.Range("J139").Value = 66238.664
$ws.Range("L139").Value = 66238.664
$ws.Range("N139").Value = -76518.664

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1103.5385
$ws.Range("I3").Value = 445.5
$ws.Range("J3").Value = 9000
$ws.Range("K3").Value = 445.5
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -331.5
$ws.Range("N3").Value = -9228

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4155.9395
$ws.Range("I31").Value = 2677.2307
$ws.Range("J31").Value = 5117.1
$ws.Range("K31").Value = 2677.2307
$ws.Range("L31").Value = 5117.1
$ws.Range("M31").Value = -2382.2307
$ws.Range("N31").Value = -5707.1
# Row 34
$ws.Range("H34").Value = 4155.9395
$ws.Range("I34").Value = 2677.2307
$ws.Range("J34").Value = 5117.1
$ws.Range("K34").Value = 2677.2307
$ws.Range("L34").Value = 5117.1
$ws.Range("M34").Value = -2475.2307
$ws.Range("N34").Value = -5521.1
# Row 140
$ws.Range("H140").Value = 119000
$ws.Range("J140").Value = 119000
$ws.Range("L140").Value = 119000
$ws.Range("N140").Value = -129360

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1974547.2
$ws.Range("J113").Value = 2333192.2
$ws.Range("L113").Value = 6999576.600000001
$ws.Range("N113").Value = -7003916.600000001
# Row 122
$ws.Range("H122").Value = 1082896
$ws.Range("I122").Value = 703.63635
$ws.Range("K122").Value = 6332.72715
$ws.Range("M122").Value = -3882.72715
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
# Row 129
$ws.Range("H129").Value = 1748.5555
$ws.Range("J129").Value = 2663
$ws.Range("L129").Value = 7989
$ws.Range("N129").Value = -17989
# Row 131
$ws.Range("H131").Value = 21161.545
$ws.Range("I131").Value = 2500.75
$ws.Range("K131").Value = 7502.25
$ws.Range("M131").Value = -2462.25

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2035.6364
$ws.Range("I102").Value = 913.96295
$ws.Range("K102").Value = 913.96295
$ws.Range("M102").Value = 708.03705
# Row 107
$ws.Range("H107").Value = 678.55554
$ws.Range("J107").Value = 638
$ws.Range("L107").Value = 638
$ws.Range("N107").Value = -4478
# Row 122
$ws.Range("H122").Value = 507760.62
$ws.Range("I122").Value = 696671
$ws.Range("J122").Value = 3999.6667
$ws.Range("K122").Value = 2090013
$ws.Range("L122").Value = 11999.0001
$ws.Range("M122").Value = -2087563
$ws.Range("N122").Value = -16899.0001
# Row 126
$ws.Range("H126").Value = 7368.1665
$ws.Range("I126").Value = 9934.714
$ws.Range("K126").Value = 29804.142
$ws.Range("M126").Value = -27334.142
# Row 137
$ws.Range("H137").Value = 97611.336
$ws.Range("J137").Value = 97611.336
$ws.Range("L137").Value = 97611.336
$ws.Range("N137").Value = -107811.336

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 71432500
$ws.Range("I7").Value = 166669330
$ws.Range("J7").Value = 4875
$ws.Range("K7").Value = 166669330
$ws.Range("L7").Value = 4875
$ws.Range("M7").Value = -166669218
$ws.Range("N7").Value = -5099
# Row 40
$ws.Range("H40").Value = 4693.227
$ws.Range("I40").Value = 3877.3
$ws.Range("K40").Value = 3877.3
$ws.Range("M40").Value = -3741.3
# Row 46
$ws.Range("H46").Value = 3526.2307
$ws.Range("I46").Value = 1410.5
$ws.Range("J46").Value = 3910.9092
$ws.Range("K46").Value = 1410.5
$ws.Range("L46").Value = 3910.9092
$ws.Range("M46").Value = -1222.5
$ws.Range("N46").Value = -4286.9092
# Row 61
$ws.Range("H61").Value = 1793.6471
$ws.Range("I61").Value = 1593.3125
$ws.Range("K61").Value = 1593.3125
$ws.Range("M61").Value = -1391.3125
# Row 113
$ws.Range("H113").Value = 1793.6471
$ws.Range("I113").Value = 1593.3125
$ws.Range("K113").Value = 1593.3125
$ws.Range("M113").Value = 576.6875
# Row 122
$ws.Range("H122").Value = 9136.182000000001
$ws.Range("I122").Value = 11800
$ws.Range("K122").Value = 35400
$ws.Range("M122").Value = -32950
# Row 126
$ws.Range("H126").Value = 71432500
$ws.Range("I126").Value = 166669330
$ws.Range("J126").Value = 4875
$ws.Range("K126").Value = 500007990
$ws.Range("L126").Value = 14625
$ws.Range("M126").Value = -500005520
$ws.Range("N126").Value = -19565
# Row 137
$ws.Range("H137").Value = 55552
$ws.Range("I137").Value = 54000
$ws.Range("J137").Value = 57104
$ws.Range("K137").Value = 54000
$ws.Range("L137").Value = 57104
$ws.Range("M137").Value = -48900
$ws.Range("N137").Value = -67304

$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 221999.6
$ws.Range("I43").Value = 16666.334
$ws.Range("K43").Value = 16666.334
$ws.Range("M43").Value = -16517.334
# Row 58
$ws.Range("H58").Value = 27000
$ws.Range("J58").Value = 27000
$ws.Range("L58").Value = 27000
$ws.Range("N58").Value = -27616
# Row 64
$ws.Range("H64").Value = 14979
$ws.Range("I64").Value = 14979
$ws.Range("K64").Value = 14979
$ws.Range("M64").Value = -14731
# Row 67
$ws.Range("H67").Value = 14979
$ws.Range("I67").Value = 14979
$ws.Range("K67").Value = 14979
$ws.Range("M67").Value = -14121
# Row 135
$ws.Range("H135").Value = 39499.668
$ws.Range("J135").Value = 39499.668
$ws.Range("L135").Value = 39499.668
$ws.Range("N135").Value = -49639.668
